$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'56.369.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +9.80%  "
$ws.Range("D3").Value = "'3.230.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.27%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'397.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.51%  "
$ws.Range("D6").Value = "'111.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.63%  "
$ws.Range("E7").Value = "  +2.96%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.621"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.34%  "
$ws.Range("D10").Value = "'39.35"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.57%  "
$ws.Range("D11").Value = "'0.0914"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.74%  "
$ws.Range("E12").Value = "  +2.08%  "
$ws.Range("D13").Value = "'3.735.58"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.37%  "
$ws.Range("D14").Value = "'8.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.12%  "
$ws.Range("D15").Value = "'19.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.22%  "
$ws.Range("D16").Value = "'3.231.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.41%  "
$ws.Range("E17").Value = "  +5.05%  "
$ws.Range("D18").Value = "'10.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.69%  "
$ws.Range("D19").Value = "'56.184.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +9.26%  "
$ws.Range("E20").Value = "  +2.89%  "
$ws.Range("E21").Value = "  +6.62%  "
$ws.Range("D22").Value = "'13.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.61%  "
$ws.Range("D23").Value = "'299.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +12.60%  "
$ws.Range("D24").Value = "'75.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.42%  "
$ws.Range("D25").Value = "'3.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.85%  "
$ws.Range("D26").Value = "'8.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.38%  "
$ws.Range("D27").Value = "'28.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.72%  "
$ws.Range("E28").Value = "  +3.27%  "
$ws.Range("E29").Value = "  +4.37%  "
$ws.Range("E30").Value = "  +0.38%  "
$ws.Range("D31").Value = "'0.111"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.59%  "
$ws.Range("D32").Value = "'11.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.83%  "
$ws.Range("D33").Value = "'0.0493"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.22%  "
$ws.Range("D34").Value = "'36.66"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.45%  "
$ws.Range("E35").Value = "  +3.11%  "
$ws.Range("D36").Value = "'51.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.15%  "
$ws.Range("D37").Value = "'3.13"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +25.33%  "
$ws.Range("D38").Value = "'3.53"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.35%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").Value = "'136.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.76%  "
$ws.Range("D41").Value = "'17.46"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.34%  "
$ws.Range("E42").Value = "  +3.19%  "
$ws.Range("D43").Value = "'4.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.96%  "
$ws.Range("E44").Value = "  +3.22%  "
$ws.Range("E45").Value = "  -1.80%  "
$ws.Range("D46").Value = "'22.23"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("E47").Value = "  +54.42%  "
$ws.Range("E48").Value = "  -1.99%  "
$ws.Range("D49").Value = "'2.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("D50").Value = "'2.121.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.23%  "
$ws.Range("D51").Value = "'0.0361"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.60%  "
